$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.226.74"
$ws.Range("E2").Value = "  -0.99%  "
$ws.Range("D3").Value = "1.783.74"
$ws.Range("E3").Value = "  -2.58%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "224.17"
$ws.Range("E5").Value = "  -2.92%  "
$ws.Range("E6").Value = "  +0.19%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "32.03"
$ws.Range("E8").Value = "  +0.86%  "
$ws.Range("E9").Value = "  -1.56%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0658"
$ws.Range("E10").Value = "  -2.51%  "
$ws.Range("E11").Value = "  -0.02%  "
$ws.Range("D12").Value = "2.039.35"
$ws.Range("E12").Value = "  -2.66%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.18"
$ws.Range("E13").Value = "  +7.03%  "
$ws.Range("D14").Value = "1.773.75"
$ws.Range("E14").Value = "  -3.14%  "
$ws.Range("E15").Value = "  -3.89%  "
$ws.Range("D16").Value = "34.214.35"
$ws.Range("E16").Value = "  -0.87%  "
$ws.Range("E17").Value = "  -1.38%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.78"
$ws.Range("E18").Value = "  -1.50%  "
$ws.Range("E19").Value = "  -2.41%  "
$ws.Range("E20").Value = "  -1.93%  "
$ws.Range("E21").Value = "  +0.08%  "
$ws.Range("E22").Value = "  -2.20%  "
$ws.Range("E23").Value = "  -3.98%  "
$ws.Range("E24").Value = "  -3.41%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.24"
$ws.Range("E25").Value = "  -0.61%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "16.39"
$ws.Range("E26").Value = "  -2.19%  "
$ws.Range("E27").Value = "  -1.94%  "
$ws.Range("E28").Value = "  -1.40%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.21%  "
$ws.Range("E30").Value = "  -3.33%  "
$ws.Range("E31").Value = "  -1.49%  "
$ws.Range("E32").Value = "  -2.11%  "
$ws.Range("E33").Value = "  -0.08%  "
$ws.Range("E34").Value = "  +3.42%  "
$ws.Range("D35").Value = "1.441.56"
$ws.Range("E35").Value = "  -7.37%  "
$ws.Range("E36").Value = "  -2.63%  "
$ws.Range("E37").Value = "  -1.67%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.624"
$ws.Range("E38").Value = "  -2.23%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.84"
$ws.Range("E39").Value = "  +0.93%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "82.86"
$ws.Range("E40").Value = "  -2.78%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.35"
$ws.Range("E41").Value = "  +0.22%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.888"
$ws.Range("E42").Value = "  -3.43%  "
$ws.Range("E43").Value = "  -5.15%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0509"
$ws.Range("E44").Value = "  -3.07%  "
$ws.Range("E45").Value = "  -2.10%  "
$ws.Range("E46").Value = "  +0.31%  "
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.26"
$ws.Range("E47").Value = "  -1.40%  "
$ws.Range("B48").Value = "RocketPoolETH"
$ws.Range("C48").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D48").Value = "1.939.77"
$ws.Range("E48").Value = "  -2.36%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.999"
$ws.Range("E49").Value = "  +0.04%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "98.32"
$ws.Range("E50").Value = "  +0.12%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "49.29"
$ws.Range("E51").Value = "  -6.86%  "
